# Fix ADS8686 default values, all must be hex values of the form 0xVALUE
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADS8686")

# Rows whose "Default Value" (column C) must become the hex string "0x00"
$zeroRows = @(2,3,11,13,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44)
foreach ($r in $zeroRows) {
    $ws.Cells.Item($r, 3).Value = "0x00"
}

# Row 12's Default Value was a plain 2, should become the hex string "0x02"
$ws.Cells.Item(12, 3).Value = "0x02"
